{"js": "// Merge the old title block's first two paragraphs:\n//   [italic]House of Hospitality[/italic], <br>Chapter Eleven =========================\n//   [bold]By Dorothy Day[/bold]\n// into a single, unformatted paragraph reading \"% Dorothy Day\"\n// (a pandoc-style title-block author line).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstPara = paragraphs.items[0]; // \"House of Hospitality, \\nChapter Eleven =========================\"\nconst secondPara = paragraphs.items[1]; // \"By Dorothy Day\" (bold)\nconst thirdPara = paragraphs.items[2]; // \"1939, Chapter Eleven, pp. 191 - 204.\"\n\n// Insert a brand-new, plain (unformatted) paragraph right before the third\n// paragraph; inserting fresh text this way \u2014 rather than editing the old\n// runs in place \u2014 avoids carrying over the old italic/bold run formatting.\nthirdPara.insertParagraph(\"% Dorothy Day\", \"Before\");\n\n// Drop the old two-paragraph title block entirely.\nfirstPara.delete();\nsecondPara.delete();\n\nawait context.sync();\n", "ps1": "# Merge the old title block's first two paragraphs:\n#   [italic]House of Hospitality[/italic], <br>Chapter Eleven =========================\n#   [bold]By Dorothy Day[/bold]\n# into a single, unformatted paragraph reading \"% Dorothy Day\"\n# (a pandoc-style title-block author line).\n\n$d = $word.ActiveDocument\n\n# Paragraph 3 is \"1939, Chapter Eleven, pp. 191 - 204.\" \u2014 insert a brand-new,\n# plain (unformatted) paragraph right before it. Inserting fresh text this\n# way \u2014 rather than editing the old runs in place \u2014 avoids carrying over the\n# old italic/bold run formatting.\n$thirdPara = $d.Paragraphs(3).Range\n$thirdPara.InsertParagraphBefore()\n\n$newPara = $d.Paragraphs(3).Range\n$newPara.Text = \"% Dorothy Day\"\n\n# Drop the old two-paragraph title block entirely (paragraphs 1 and 2).\n$d.Paragraphs(1).Range.Delete()\n$d.Paragraphs(1).Range.Delete()\n"}
